$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Patient identification (1er Apellido / 2do Apellido / 1er Nombre / 2do Nombre)
$ws.Range("A6").Value = "BONILLA    RUMILDA  "
# No. Expediente Clinico
$ws.Range("G6").Value = "/201770907"

# Fecha de Nacimiento - keep as text (leading quote keeps General format but forces text)
$ws.Range("A9").Formula = "'1976-07-05"
# Edad - keep as text
$ws.Range("D9").Formula = "'41"
# Lugar de Nacimiento
$ws.Range("E9").Value = "TECPAN GUATEMALA, CHIMALTENANGO"

# Estado Civil
$ws.Range("A11").Value = "CASADO"

# Documento de Identificacion
$ws.Range("G11").Value = "NO PRESENTO"

# Emergency contact: Nombre / Parentesco / Direccion / Telefono
$ws.Range("A13").Value = "JUVENTINO UMUL"
$ws.Range("D13").Value = "ESPÒSO"
$ws.Range("E13").Value = "30 CALLE 5-40 ZONA 10 MIXCO LA COMUNIDAD ALBOR"
# Telefono - keep as text
$ws.Range("G13").Formula = "'48975553"

# Fecha de la asistencia Medica: Hora / Area de urgencia / Fecha
$ws.Range("D14").Value = "Hora: 7:13:54"
$ws.Range("E14").Value = "Area de urgencia: GINECOLOGIA"
$ws.Range("A15").Value = "14/11/2017"

# Clear "Tipo de Consulta" value (consulta_externa removed)
$ws.Range("D16").Value = ""
